$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell H1 "Save" - copy formatting (bold/border/alignment) from G1 header
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add new data cell H2 with numeric value 0
$ws.Range("H2").Value = 0
